$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, shifting existing data rows (and their formatting) down by one.
$ws.Rows("2:2").Insert()

# The inserted row inherits formatting from the row above (the header); instead, match the
# style used by the data rows by copying it from the row immediately below (old row 2, now row 3).
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Force column E to keep its text representation instead of being reinterpreted as a date.
$ws.Cells.Item(2, 5).NumberFormat = "@"

# Populate the new top data row with the latest circular entry.
$ws.Cells.Item(2, 1).Value = 10
$ws.Cells.Item(2, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(2, 3).Value = "IE07"
$ws.Cells.Item(2, 4).Value = 277.95
$ws.Cells.Item(2, 5).Value = "01-10-2025"
$ws.Cells.Item(2, 6).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf"

$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf")
